$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.043.05"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").Value = "'1.689.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "'216.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'23.75"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.77%  "

$ws.Range("D9").Value = "'0.264"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.65%  "

$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "'1.928.86"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").Value = "'1.686.34"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("E14").Value = "  +0.14%  "

$ws.Range("D15").Value = "'0.557"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").Value = "'67.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "'250.84"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.52%  "

$ws.Range("D18").Value = "'28.058.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.53%  "

$ws.Range("D19").Value = "'0.0₃0742"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.12%  "

$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "'9.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("D25").Value = "'147.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").Value = "'7.33"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.88%  "

$ws.Range("D27").Value = "'16.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.02%  "

$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("E30").Value = "  +6.70%  "

$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").Value = "'3.38"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("D33").Value = "'3.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("D34").Value = "'1.426.24"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.52%  "

$ws.Range("D35").Value = "'1.61"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.31%  "

$ws.Range("D36").Value = "'0.939"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "'2.38"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("D38").Value = "'0.591"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.30%  "

$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("E40").Value = "  -3.31%  "

$ws.Range("D41").Value = "'69.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").Value = "'5.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.96%  "

$ws.Range("D44").Value = "'1.836.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "'0.800"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.94%  "

$ws.Range("D47").Value = "'1.71"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.77%  "

$ws.Range("D48").Value = "'89.46"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  -1.06%  "

$ws.Range("D51").Value = "'7.81"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.10%  "
